# Apply "Penalty Reward System" edits (unfinished) to the PO data workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update existing rows 10-12 with new date/quantity values
$ws1.Cells.Item(10, 2).Value = 60

$ws1.Cells.Item(11, 1).Value = 45109.99999999999
$ws1.Cells.Item(11, 2).Value = 140

$ws1.Cells.Item(12, 1).Value = 45116.99999999999
$ws1.Cells.Item(12, 2).Value = 40

# Remove the now-obsolete trailing rows (13-16)
$ws1.Rows("13:16").Delete() | Out-Null

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(6, 2).Value = 200
$ws2.Cells.Item(7, 2).Value = 40
